# Adds the "ITERACOES 2000" worksheet (planilha de 2000 iteracoes do fss),
# modeled after the existing "ITERACOES 1000" sheet: same layout/labels,
# new data values for 2000 iterations.

$wb = $excel.ActiveWorkbook

# "ITERACOES 1000" is the 3rd tab in this workbook.
$src = $wb.Worksheets.Item(3)

# Duplicate it, placing the copy right after the source sheet (matches how
# Excel's own "Move or Copy..." -> Create a copy works), then rename it.
$src.Copy($null, $src)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "ITERACOES 2000"
$ws.Activate()

# Row 4 label: "1000 iterações" -> "2000 iterações"
$ws.Range("A4").Value2 = "2000 iterações"

# New per-run data (rows 10-39) and totals (row 40) for columns B:E
# (B=10 peixes, C=20 peixes, D=50 peixes, E=100 peixes)
$data = @(
    @(0.442496084537356,      0.356542546615636,     0.41978974555784898, 0.59193112215383203),
    @(0.43086871113966801,    0.42956881749951598,   0.422502271960499,   0.50518965934307303),
    @(0.415997502782329,      0.43409864389493802,   0.44172138289287599, 0.50050155177436795),
    @(0.4075272468958,        0.43501254349705698,   0.44424541467573803, 0.49650075441082298),
    @(0.41756590683755301,    0.42652647000835903,   0.45241504807212901, 0.49457315828410198),
    @(0.41858404718127701,    0.42806980673726902,   0.45715753747021898, 0.49067558730829702),
    @(0.41712659420536502,    0.43039761873034799,   0.460998246260163,   0.49261714149757602),
    @(0.41165725472041598,    0.43494654147317402,   0.46542657696082701, 0.49526489718225403),
    @(0.41382105990315698,    0.44104029266287897,   0.46498828213769899, 0.49076969474724103),
    @(0.41591321216376897,    0.44417045263841798,   0.46535283076218398, 0.49305166146585599),
    @(0.42049371198825503,    0.44569323298236402,   0.46669218346089603, 0.49604504729885401),
    @(0.42070268068577299,    0.446651136427138,     0.463241733953281,   0.49558805103791198),
    @(0.42198934465319898,    0.447908045723589,     0.46606882685753098, 0.49410009911707198),
    @(0.42404327417778698,    0.44566897636008801,   0.46476799228559101, 0.49867633241752801),
    @(0.427648435083678,      0.451009137401918,     0.46748873490668902, 0.49899940824634198),
    @(0.42611285261431397,    0.45277101716294699,   0.46960966995350401, 0.50105902268361702),
    @(0.42688246462634899,    0.45139470905720303,   0.47141518180715902, 0.50082237207601199),
    @(0.42788661627768698,    0.45253420980249298,   0.471745678765773,   0.49911267879642601),
    @(0.42859435539555901,    0.452711957497622,     0.47435304020555102, 0.49973851056142699),
    @(0.42674959990081501,    0.45304094033761599,   0.47525896523388,    0.50185622995266999),
    @(0.42789267357517402,    0.454153658845104,     0.47488783807352702, 0.501870216262015),
    @(0.42658623458348099,    0.457030123818565,     0.47455656195282597, 0.50114247739387097),
    @(0.42752780295738502,    0.45815744739179998,   0.475186100335935,   0.50313361764063702),
    @(0.42771234980212602,    0.45732264107385501,   0.47520459656550501, 0.50282887834676204),
    @(0.42564939725103901,    0.45656920734150702,   0.47471837581683901, 0.50173612169128001),
    @(0.42526094362174999,    0.45586345820069502,   0.47587054275225399, 0.501522292936325),
    @(0.426246757103801,      0.456792889563295,     0.47645153627650499, 0.50111257933980502),
    @(0.426615536590286,      0.45776236992435398,   0.476047379301481,   0.50140426171459795),
    @(0.42764992296110599,    0.45783340614117701,   0.47644290783953902, 0.50108045561947701),
    @(0.42933537131421001,    0.458352593263265,     0.47564473431758297, 0.50060582681427201)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 10 + $i
    $ws.Cells.Item($row, 2).Value2 = $data[$i][0]
    $ws.Cells.Item($row, 3).Value2 = $data[$i][1]
    $ws.Cells.Item($row, 4).Value2 = $data[$i][2]
    $ws.Cells.Item($row, 5).Value2 = $data[$i][3]
}

# Totals row
$ws.Cells.Item(40, 2).Value2 = 8989642
$ws.Cells.Item(40, 3).Value2 = 20692449
$ws.Cells.Item(40, 4).Value2 = 49221669
$ws.Cells.Item(40, 5).Value2 = 98042792

# The pasted-in values for column B (10 peixes) and column D (50 peixes)
# come in as plain, un-styled numbers rather than carrying over the old
# sheet's per-cell formatting, so drop back to the default style there.
$ws.Range("B10:B39").Style = "Normal"
$ws.Range("D10:D40").Style = "Normal"

# Scroll/selection state on the new tab.
$ws.Range("F26").Select()

Write-Host "Added ITERACOES 2000 sheet with" $wb.Worksheets.Count "total sheets"
